$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 3) with the new item's information
$ws.Range("A3").Value = "Furniture, Fixtures and Equipment-CPGC"
$ws.Range("B3").Value = "FFE-FUR-CPGC-1030"
$ws.Range("D3").Value = "ste"
$ws.Range("I3").Value = "drum/s"
$ws.Range("K3").Value = "Jan Lester Mercene Madriaga"
$ws.Range("M3").Value = "Operations"

# Re-fit the columns so the new (longer/shorter) content is fully visible,
# matching the widths Excel's own "best fit" produced for this content.
$ws.Columns.Item(1).ColumnWidth = 45
$ws.Columns.Item(2).ColumnWidth = 20.333333333333332
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(9).ColumnWidth = 7.333333333333333
$ws.Columns.Item(11).ColumnWidth = 32.166666666666664
